$d = $word.ActiveDocument

# 1) Main body: the bold "TERE" placeholder (right after "A ") -> "QWER"
$bodyRng = $d.Content
$bodyRng.Find.Execute("TERE", $true, $true, $false, $false, $false, $true, 0, $false, "QWER", 1) | Out-Null

# 2) Header placeholders, in left-to-right document order. Several runs share
#    the same old text ("Tre" / "tre") but map to different new text, so each
#    occurrence is located and replaced one at a time (Find auto-advances the
#    range forward after every successful match).
$sec = $d.Sections(1)
$hdr = $sec.Headers(1)
$rng = $hdr.Range.Duplicate

$pairs = @(
    @("TRE", "QWER"),
    @("TERE", "QWER"),
    @("Tre", "Qwer"),
    @("Tre", "Qwer"),
    @("Tre", "Qewr"),
    @("Tre", "Qewr"),
    @("Tre", "Qwer"),
    @("tre", "qwer"),
    @("tre", "qwer"),
    @("tre", "qwer")
)

foreach ($pair in $pairs) {
    $needle = $pair[0]
    $replacement = $pair[1]
    $rng.Find.Execute($needle, $true, $true, $false, $false, $false, $true, 0, $false, $replacement, 1) | Out-Null
}
